$wb = $excel.ActiveWorkbook

# Add the new worksheet after the existing one
$existing = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $existing)
$ws.Name = "labor_beneimp_stochastic_reg"

# Row 3 headers (same values/order as sheet1)
$headers = @("mean", "p10", "p25", "p50", "p75", "p90", "p99")

$col = 2  # column B
foreach ($h in $headers) {
    $ws.Cells.Item(3, $col).Value = $h
    $col++
}
foreach ($h in $headers) {
    $ws.Cells.Item(3, $col).Value = $h
    $col++
}

# Row 4 data (A4 must stay text "2019", matching the shared string used in sheet1!A4)
# Build the text value via a formula in a scratch cell, then paste-special
# (values only) into A4 so the cell becomes a plain text/shared-string cell
# without picking up any number-format style override.
$ws.Cells.Item(1, 1).Formula = "=""2019"""
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(4, 1).PasteSpecial(-4163)
$ws.Cells.Item(1, 1).ClearContents()
$excel.CutCopyMode = $false

$ws.Cells.Item(4, 2).Value = 863689.71242246288
$ws.Cells.Item(4, 3).Value = 55419.078125
$ws.Cells.Item(4, 4).Value = 200000
$ws.Cells.Item(4, 5).Value = 302285.875
$ws.Cells.Item(4, 6).Value = 860511.875
$ws.Cells.Item(4, 7).Value = 2000000
$ws.Cells.Item(4, 8).Value = 9505159

$ws.Cells.Item(4, 9).Value = 912108.76194957457
$ws.Cells.Item(4, 10).Value = 60457.17578125
$ws.Cells.Item(4, 11).Value = 200000
$ws.Cells.Item(4, 12).Value = 340603.6875
$ws.Cells.Item(4, 13).Value = 868189.0625
$ws.Cells.Item(4, 14).Value = 2015239.25
$ws.Cells.Item(4, 15).Value = 9505159

$excel.ActiveWorkbook.Worksheets.Item(1).Activate()
